# Regenerate the "K" column (column G) values for the save_data sheet.
# The scraping/regeneration pipeline now derives K (strikeouts) differently
# ("use K instead of Strike#"), so the previously stored values in column G
# are replaced with the newly computed ones for each data row (rows 2-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = [ordered]@{
    2  = 1
    3  = 1
    4  = 0
    5  = 0
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 2
    11 = 0
    12 = 2
    13 = 1
    14 = 2
    15 = 1
    16 = 0
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 0
    22 = 1
    23 = 1
    24 = 2
    25 = 2
    26 = 1
    27 = 0
    28 = 0
    29 = 1
    30 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
